$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) column F for the two affected rows
# in both the "展览" sheet and the duplicated "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 120
    $ws.Range("F4").Value = 68
}

$wb.Save()
